$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of locations (rows 2-17, column A), ranks in column B stay 1-16.
$locations = @(
    "MAŁOPOLSKIE",
    "POMORSKIE",
    "DOLNOŚLĄSKIE",
    "ŚLĄSKIE",
    "MAZOWIECKIE",
    "ŁÓDZKIE",
    "ZACHODNIOPOMORSKIE",
    "PODKARPACKIE",
    "LUBUSKIE",
    "LUBELSKIE",
    "ŚWIĘTOKRZYSKIE",
    "WARMIŃSKO-MAZURSKIE",
    "PODLASKIE",
    "OPOLSKIE",
    "WIELKOPOLSKIE",
    "KUJAWSKO-POMORSKIE"
)

for ($i = 0; $i -lt $locations.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $locations[$i]
}
